$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 6, shifting existing rows 6-13 down to 8-15
$ws.Rows.Item(6).Resize(2).Insert()

$ws.Range("A6").Value = 'fgWKXlPr'
$ws.Range("B6").Value = '21/11/2024'
$ws.Range("C6").Value = '20:00'
$ws.Range("D6").Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Range("E6").Value = 'Tomayapo'
$ws.Range("F6").Value = 'Bolivar'
$ws.Range("G6").Value = 4.1
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 1.75
$ws.Range("J6").Value = 4.5
$ws.Range("K6").Value = 2.3
$ws.Range("L6").Value = 2.3
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 15
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.33
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 2.15
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 13
$ws.Range("X6").Value = 23
$ws.Range("Y6").Value = 13
$ws.Range("Z6").Value = 41
$ws.Range("AA6").Value = 29
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 15
$ws.Range("AD6").Value = 7.5
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 41
$ws.Range("AG6").Value = 8.5
$ws.Range("AH6").Value = 9.5
$ws.Range("AI6").Value = 8.5
$ws.Range("AJ6").Value = 15
$ws.Range("AK6").Value = 13
$ws.Range("AL6").Value = 21
$ws.Range("AM6").Value = 151
$ws.Range("AN6").Value = 6
$ws.Range("AO6").Value = 21
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 67
$ws.Range("AR6").Value = 81
$ws.Range("AS6").Value = 151
$ws.Range("AT6").Value = 3.25
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 4
$ws.Range("AX6").Value = 9
$ws.Range("AY6").Value = 17
$ws.Range("AZ6").Value = 29
$ws.Range("BA6").Value = 41
$ws.Range("BB6").Value = 101
$ws.Range("A7").Value = 'zPs4uVuR'
$ws.Range("B7").Value = '21/11/2024'
$ws.Range("C7").Value = '20:00'
$ws.Range("D7").Value = 'BRAZIL - SERIE A BETANO'
$ws.Range("E7").Value = 'Vasco'
$ws.Range("F7").Value = 'Internacional'
$ws.Range("G7").Value = 3.6
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.2
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 2.88
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 9.5
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 41
$ws.Range("AA7").Value = 29
$ws.Range("AB7").Value = 41
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 7
$ws.Range("AH7").Value = 9.5
$ws.Range("AI7").Value = 9.5
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 19
$ws.Range("AL7").Value = 34
$ws.Range("AM7").Value = 351
$ws.Range("AN7").Value = 5.5
$ws.Range("AO7").Value = 21
$ws.Range("AP7").Value = 29
$ws.Range("AQ7").Value = 67
$ws.Range("AR7").Value = 101
$ws.Range("AS7").Value = 251
$ws.Range("AT7").Value = 2.5
$ws.Range("AU7").Value = 8.5
$ws.Range("AV7").Value = 67
$ws.Range("AW7").Value = 4
$ws.Range("AX7").Value = 12
$ws.Range("AY7").Value = 23
$ws.Range("AZ7").Value = 41
$ws.Range("BA7").Value = 67
$ws.Range("BB7").Value = 201
$ws.Range("BC7").Value = 126
$ws.Range("BD7").Value = 126
